$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (prices + volume deltas).
# Columns B-E are plain text in the source sheet; D-column numeric-looking
# price strings need an explicit Text format so COM does not coerce them
# into floating point numbers (which would lose trailing zeros / precision).

$ws.Range("D2").Value = "23.057.32"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").Value = "1.598.70"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.25"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3777"
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3636"
$ws.Range("E8").Value = "  -0.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.61"
$ws.Range("E9").Value = "  +1.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.249"
$ws.Range("E10").Value = "  -2.25%  "

$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08130"
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.31"
$ws.Range("E13").Value = "  -2.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.568"
$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.353"
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001243"
$ws.Range("E16").Value = "  -2.08%  "

$ws.Range("D17").Value = "1.598.96"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.82"
$ws.Range("E18").Value = "  +0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06833"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.13"
$ws.Range("E20").Value = "  -2.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.511"
$ws.Range("E21").Value = "  -1.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.99"
$ws.Range("E23").Value = "  -1.47%  "

$ws.Range("D24").Value = "23.061.77"
$ws.Range("E24").Value = "  -0.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.368"
$ws.Range("E25").Value = "  -1.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.769"
$ws.Range("E26").Value = "  -6.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.05"
$ws.Range("E27").Value = "  -0.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.95"
$ws.Range("E28").Value = "  -1.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.256"
$ws.Range("E29").Value = "  -1.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.21"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.353"
$ws.Range("E31").Value = "  -4.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.794"
$ws.Range("E32").Value = "  -8.33%  "

$ws.Range("D33").Value = "1.774.49"
$ws.Range("E33").Value = "  +0.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9574"
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07549"
$ws.Range("E35").Value = "  -2.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.15"
$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.190"
$ws.Range("E37").Value = "  -2.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02695"
$ws.Range("E38").Value = "  -3.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2510"
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08815"
$ws.Range("E40").Value = "  -0.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.355"
$ws.Range("E41").Value = "  -2.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7013"
$ws.Range("E42").Value = "  -2.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.26"
$ws.Range("E43").Value = "  -4.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.21"
$ws.Range("E44").Value = "  -5.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6560"
$ws.Range("E45").Value = "  -1.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.992"
$ws.Range("E47").Value = "  +0.14%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.271"
$ws.Range("E48").Value = "  -2.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.58"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07939"
$ws.Range("E50").Value = "  -1.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.214"
$ws.Range("E51").Value = "  +3.12%  "
